$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: "模型名称" -> "Lambda" (E2) ---
$ws.Range("E2").Value = "Lambda"

# --- Row 22: add F22 = 0.1 ---
$ws.Range("F22").Value = 0.1

# --- Row 27: add date, add E27 comment, remove N27 (value moved out) ---
$ws.Range("B27").Value = 43321
$ws.Range("B27").NumberFormat = "MM/DD/YY"
$ws.Range("E27").Value = "loc=2, kl=0.5"
$ws.Range("N27").ClearContents()

# --- Row 28: add date, add E28 comment, remove N28 ---
$ws.Range("B28").Value = 43321
$ws.Range("B28").NumberFormat = "MM/DD/YY"
$ws.Range("E28").Value = "loc=1, kl=0.3"
$ws.Range("N28").ClearContents()

# --- New row 30 ---
$ws.Range("A30").Value = "3-1"
$ws.Range("C30").Value = "bug when concate sketchInput"

# --- New row 31 ---
$ws.Range("A31").Value = "3-2"
$ws.Range("C31").Value = "bug"
$ws.Range("E31").Value = "loc=0.1,kl=0.2"

# --- New row 32 ---
$ws.Range("A32").Value = "3-3"
$ws.Range("C32").Value = "bug"
$ws.Range("E32").Value = "eof=0.1, kl=0"

# --- New row 33 ---
$ws.Range("A33").Value = "3-4"
$ws.Range("B33").Value = 43322
$ws.Range("B33").NumberFormat = "MM/DD/YY"
$ws.Range("E33").Value = "Kl=0"
$ws.Range("G33").Value = "Recons_High, KL_High"
$ws.Range("J33").Value = 60000
$ws.Range("L33").Value = "[40000, 55000]"

# --- New row 34 ---
$ws.Range("A34").Value = "3-5"
$ws.Range("E34").Value = "Kl=0.01"

# --- Column widths: split column E out of A:E group, split column L out of K:L group ---
$ws.Columns.Item(5).ColumnWidth = 17.1
$ws.Columns.Item(12).ColumnWidth = 5.6

# --- Update selection to A35 (matches final cursor position in the authored file) ---
$ws.Range("A35").Select()
